# Swap the match data between row 117 and row 118 on the "Lithuania A Lyga"
# sheet. Columns A (row index), C (Div), D (Div Original Name) and E (Date)
# are identical for both rows and stay untouched; column B (id) and the
# F:AC block (HomeTeam .. PL_AhUnder) are exchanged between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Lithuania A Lyga")

$row1 = 117
$row2 = 118

# Numeric column indices: B=2, F=6 .. AC=29
$colIndices = @(2) + @(6..29)

foreach ($col in $colIndices) {
    $cell1 = $ws.Cells.Item($row1, $col)
    $cell2 = $ws.Cells.Item($row2, $col)

    $val1 = $cell1.Value2
    $val2 = $cell2.Value2

    $cell1.Value2 = $val2
    $cell2.Value2 = $val1
}
